$d = $word.ActiveDocument

# Update the date line (first paragraph)
$d.Paragraphs.Item(1).Range.Text = "2026-01-04 Sunday"

# Update the data rows of the table
$t = $d.Tables.Item(1)

# Row 1
$t.Cell(1, 1).Range.Text = "20÷7=2, 6"
$t.Cell(1, 2).Range.Text = "83÷5=16, 3"
$t.Cell(1, 3).Range.Text = "74÷6=12, 2"
$t.Cell(1, 4).Range.Text = "67÷5=13, 2"
$t.Cell(1, 5).Range.Text = "78÷6=13, 0"

# Row 5
$t.Cell(5, 1).Range.Text = "33÷2=16, 1"
$t.Cell(5, 2).Range.Text = "34÷8=4, 2"
$t.Cell(5, 3).Range.Text = "45÷8=5, 5"
$t.Cell(5, 4).Range.Text = "34÷8=4, 2"
$t.Cell(5, 5).Range.Text = "95÷9=10, 5"

# Row 9
$t.Cell(9, 1).Range.Text = "39÷4=9, 3"
$t.Cell(9, 2).Range.Text = "17÷3=5, 2"
$t.Cell(9, 3).Range.Text = "35÷5=7, 0"
$t.Cell(9, 4).Range.Text = "20÷4=5, 0"
$t.Cell(9, 5).Range.Text = "48÷9=5, 3"

# Row 13
$t.Cell(13, 1).Range.Text = "95÷4=23, 3"
$t.Cell(13, 2).Range.Text = "18÷3=6, 0"
$t.Cell(13, 3).Range.Text = "53÷3=17, 2"
$t.Cell(13, 4).Range.Text = "18÷4=4, 2"
$t.Cell(13, 5).Range.Text = "60÷5=12, 0"

# Row 17
$t.Cell(17, 1).Range.Text = "80÷9=8, 8"
$t.Cell(17, 2).Range.Text = "25÷3=8, 1"
$t.Cell(17, 3).Range.Text = "86÷4=21, 2"
$t.Cell(17, 4).Range.Text = "92÷5=18, 2"
$t.Cell(17, 5).Range.Text = "59÷2=29, 1"

